# Add a new "2022-Q1" worksheet (mirroring the structure of the other
# quarterly holdings sheets) right before the "总计" (totals) sheet, and
# update the "总计" sheet with the new quarter's aggregate row.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item("总计")
$template = $wb.Worksheets.Item("2021-Q3")

# --- 1. Create the new quarter sheet, just before "总计" ---------------
$new = $wb.Worksheets.Add($total)
$new.Name = "2022-Q1"

# NOTE: inserting a sheet shifts "总计" to a new tab position, which
# stales out the handle captured above (the engine resolves worksheet
# handles positionally) -- re-resolve it by name before writing to it.
$total = $wb.Worksheets.Item("总计")

# Copy the whole 13-row/8-column table (values + formatting) from an
# existing quarterly sheet so the row-index column (A) and all cell
# styles/borders match exactly; text/values are then overwritten below.
$template.Range("A1:H13").Copy($new.Range("A1:H13"))

# Header row text (set explicitly; wording drifted slightly release to
# release, e.g. "基金规模" vs "基金金额", so don't just trust the copy).
$new.Range("B1").Value = "基金代码"
$new.Range("C1").Value = "基金名称"
$new.Range("D1").Value = "基金规模"
$new.Range("E1").Value = "股票总仓位"
$new.Range("F1").Value = "仓位占比"
$new.Range("G1").Value = "持有市值(亿元)"
$new.Range("H1").Value = "仓位排名"

# --- helper: write a value as literal TEXT (not auto-converted to a
# number), matching the source data's inlineStr typing, then drop the
# leftover "@" number-format so the cell keeps the sheet's default style.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# --- 2. Fill in the 2022-Q1 fund-holding rows ---------------------------
Set-TextValue $new.Range("B2") "006682"
Set-TextValue $new.Range("C2") "景顺长城中证500指数增强"
Set-TextValue $new.Range("D2") "16.63"
Set-TextValue $new.Range("E2") "87.75"
Set-TextValue $new.Range("F2") "1.83"
Set-TextValue $new.Range("G2") "0.3043"
$new.Range("H2").Value = 10

Set-TextValue $new.Range("B3") "011997"
Set-TextValue $new.Range("C3") "景顺长城安盈回报一年持有期混合型证券投资基金A"
Set-TextValue $new.Range("D3") "5.74"
Set-TextValue $new.Range("E3") "29.73"
Set-TextValue $new.Range("F3") "1.85"
Set-TextValue $new.Range("G3") "0.1062"
$new.Range("H3").Value = 4

Set-TextValue $new.Range("B4") "460009"
Set-TextValue $new.Range("C4") "华泰柏瑞量化先行混合A"
Set-TextValue $new.Range("D4") "9.13"
Set-TextValue $new.Range("E4") "90.47"
Set-TextValue $new.Range("F4") "0.89"
Set-TextValue $new.Range("G4") "0.0813"
$new.Range("H4").Value = 7

Set-TextValue $new.Range("B5") "005632"
Set-TextValue $new.Range("C5") "鹏华量化先锋混合"
Set-TextValue $new.Range("D5") "3.10"
Set-TextValue $new.Range("E5") "92.91"
Set-TextValue $new.Range("F5") "1.36"
Set-TextValue $new.Range("G5") "0.0422"
$new.Range("H5").Value = 10

Set-TextValue $new.Range("B6") "011731"
Set-TextValue $new.Range("C6") "国投瑞银安睿混合A"
Set-TextValue $new.Range("D6") "2.58"
Set-TextValue $new.Range("E6") "43.48"
Set-TextValue $new.Range("F6") "0.72"
Set-TextValue $new.Range("G6") "0.0186"
$new.Range("H6").Value = 7

Set-TextValue $new.Range("B7") "260117"
Set-TextValue $new.Range("C7") "景顺长城支柱产业混合"
Set-TextValue $new.Range("D7") "0.24"
Set-TextValue $new.Range("E7") "76.78"
Set-TextValue $new.Range("F7") "2.98"
Set-TextValue $new.Range("G7") "0.0072"
$new.Range("H7").Value = 8

Set-TextValue $new.Range("B8") "011732"
Set-TextValue $new.Range("C8") "国投瑞银安睿混合C"
Set-TextValue $new.Range("D8") "0.95"
Set-TextValue $new.Range("E8") "43.48"
Set-TextValue $new.Range("F8") "0.72"
Set-TextValue $new.Range("G8") "0.0068"
$new.Range("H8").Value = 7

Set-TextValue $new.Range("B9") "011998"
Set-TextValue $new.Range("C9") "景顺长城安盈回报一年持有期混合型证券投资基金C"
Set-TextValue $new.Range("D9") "0.23"
Set-TextValue $new.Range("E9") "29.73"
Set-TextValue $new.Range("F9") "1.85"
Set-TextValue $new.Range("G9") "0.0043"
$new.Range("H9").Value = 4

Set-TextValue $new.Range("B10") "005260"
Set-TextValue $new.Range("C10") "银华稳健增利灵活配置混合A"
Set-TextValue $new.Range("D10") "0.32"
Set-TextValue $new.Range("E10") "91.49"
Set-TextValue $new.Range("F10") "0.76"
Set-TextValue $new.Range("G10") "0.0024"
$new.Range("H10").Value = 10

Set-TextValue $new.Range("B11") "010246"
Set-TextValue $new.Range("C11") "华泰柏瑞量化先行混合C"
Set-TextValue $new.Range("D11") "0.12"
Set-TextValue $new.Range("E11") "90.47"
Set-TextValue $new.Range("F11") "0.89"
Set-TextValue $new.Range("G11") "0.0011"
$new.Range("H11").Value = 7

Set-TextValue $new.Range("B12") "006195"
Set-TextValue $new.Range("C12") "国金量化多因子股票"
Set-TextValue $new.Range("D12") "0.09"
Set-TextValue $new.Range("E12") "80.71"
Set-TextValue $new.Range("F12") "0.93"
Set-TextValue $new.Range("G12") "0.0008"
$new.Range("H12").Value = 1

Set-TextValue $new.Range("B13") "005261"
Set-TextValue $new.Range("C13") "银华稳健增利灵活配置混合C"
Set-TextValue $new.Range("D13") "0.02"
Set-TextValue $new.Range("E13") "91.49"
Set-TextValue $new.Range("F13") "0.76"
Set-TextValue $new.Range("G13") "0.0002"
$new.Range("H13").Value = 10

# --- 3. Update "总计": push existing quarters down a row and insert the
#        new 2022-Q1 aggregate at the top of the data ---------------------

# Give the new bottom row (A6) the same row-index style as the rest of
# column A before filling it in.
$total.Range("A5").Copy()
$total.Range("A6").PasteSpecial(-4122)

$total.Range("B6").Value = $total.Range("B5").Value2
$total.Range("C6").Value = $total.Range("C5").Value2
$total.Range("D6").Value = $total.Range("D5").Value2

$total.Range("B5").Value = $total.Range("B4").Value2
$total.Range("C5").Value = $total.Range("C4").Value2
$total.Range("D5").Value = $total.Range("D4").Value2

$total.Range("B4").Value = $total.Range("B3").Value2
$total.Range("C4").Value = $total.Range("C3").Value2
$total.Range("D4").Value = $total.Range("D3").Value2

$total.Range("B3").Value = $total.Range("B2").Value2
$total.Range("C3").Value = $total.Range("C2").Value2
$total.Range("D3").Value = $total.Range("D2").Value2

$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 12
$total.Range("D2").Value = 0.58

$total.Range("A6").Value = 4
